$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit the existing "testnumber" shared string (C1, only referencer) in place
# to "checkOutAsGuest" before it gets overwritten/removed, then reuse it on
# the new rows below.
$ws.Range("C1").Value = "checkOutAsGuest"
$ws.Range("A5").Value = "checkOutAsGuest"
$ws.Range("B5").Value = "Black Bootcut Jeans"
$ws.Range("A4").Value = "searchProductTest"
$ws.Range("B4").Value = "Shoes"
$ws.Range("A6").Value = "checkOutAsGuest"
$ws.Range("B6").Value = "Dakota Indigo Washed Jeans"

# Remove column C entirely (also removes the now-unused C1/C2/C3 values)
$ws.Columns.Item(3).Delete()

$ws.Columns.Item(1).ColumnWidth = 30.85546875
$ws.Columns.Item(2).ColumnWidth = 26.5703125

$ws.Range("B4").Select()
